$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.797777
$ws.Range("H2").Value = 122.393331
$ws.Range("I2").Value = 0.2689231481273683
$ws.Range("J2").Value = 0.2689231481273683
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8694479999999999
$ws.Range("N2").Value = 2.608344
$ws.Range("O2").Value = 0.1600841558454311
$ws.Range("P2").Value = 0.1600841558454311
$ws.Range("Q2").Value = 35.471545617096
$ws.Range("R2").Value = 319.243910553864
$ws.Range("S2").Value = 0.04305033515526559
$ws.Range("T2").Value = 0.0430503351552656
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.797777
$ws.Range("H3").Value = 122.393331
$ws.Range("I3").Value = 0.2689231481273683
$ws.Range("J3").Value = 0.2689231481273683
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.177221666666667
$ws.Range("N3").Value = 9.531665
$ws.Range("O3").Value = 0.5849951330524047
$ws.Range("P3").Value = 0.5849951330524048
$ws.Range("Q3").Value = 129.623581036235
$ws.Range("R3").Value = 1166.612229326115
$ws.Range("S3").Value = 0.1573187328196413
$ws.Range("T3").Value = 0.1573187328196414
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.797777
$ws.Range("H4").Value = 122.393331
$ws.Range("I4").Value = 0.2689231481273683
$ws.Range("J4").Value = 0.2689231481273683
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.384523666666666
$ws.Range("N4").Value = 4.153570999999999
$ws.Range("O4").Value = 0.2549207111021641
$ws.Range("P4").Value = 0.2549207111021641
$ws.Range("Q4").Value = 56.485487803889
$ws.Range("R4").Value = 508.369390235001
$ws.Range("S4").Value = 0.06855408015246134
$ws.Range("T4").Value = 0.06855408015246135
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 46.219831
$ws.Range("H5").Value = 138.659493
$ws.Range("I5").Value = 0.3046632285488233
$ws.Range("J5").Value = 0.3046632285488233
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.8694479999999999
$ws.Range("N5").Value = 2.608344
$ws.Range("O5").Value = 0.1600841558454311
$ws.Range("P5").Value = 0.1600841558454311
$ws.Range("Q5").Value = 40.18573962328799
$ws.Range("R5").Value = 361.6716566095919
$ws.Range("S5").Value = 0.04877175575938204
$ws.Range("T5").Value = 0.04877175575938204
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 46.219831
$ws.Range("H6").Value = 138.659493
$ws.Range("I6").Value = 0.3046632285488233
$ws.Range("J6").Value = 0.3046632285488233
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.177221666666667
$ws.Range("N6").Value = 9.531665
$ws.Range("O6").Value = 0.5849951330524047
$ws.Range("P6").Value = 0.5849951330524048
$ws.Range("Q6").Value = 146.8506484828717
$ws.Range("R6").Value = 1321.655836345845
$ws.Range("S6").Value = 0.1782265059210941
$ws.Range("T6").Value = 0.1782265059210941
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 46.219831
$ws.Range("H7").Value = 138.659493
$ws.Range("I7").Value = 0.3046632285488233
$ws.Range("J7").Value = 0.3046632285488233
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.384523666666666
$ws.Range("N7").Value = 4.153570999999999
$ws.Range("O7").Value = 0.2549207111021641
$ws.Range("P7").Value = 0.2549207111021641
$ws.Range("Q7").Value = 63.99244988883365
$ws.Range("R7").Value = 575.932048999503
$ws.Range("S7").Value = 0.0776649668683472
$ws.Range("T7").Value = 0.0776649668683472
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 64.69033266666666
$ws.Range("H8").Value = 194.070998
$ws.Range("I8").Value = 0.4264136233238083
$ws.Range("J8").Value = 0.4264136233238083
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.8694479999999999
$ws.Range("N8").Value = 2.608344
$ws.Range("O8").Value = 0.1600841558454311
$ws.Range("P8").Value = 0.1600841558454311
$ws.Range("Q8").Value = 56.24488035636799
$ws.Range("R8").Value = 506.203923207312
$ws.Range("S8").Value = 0.0682620649307835
$ws.Range("T8").Value = 0.0682620649307835
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 64.69033266666666
$ws.Range("H9").Value = 194.070998
$ws.Range("I9").Value = 0.4264136233238083
$ws.Range("J9").Value = 0.4264136233238083
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.177221666666667
$ws.Range("N9").Value = 9.531665
$ws.Range("O9").Value = 0.5849951330524047
$ws.Range("P9").Value = 0.5849951330524048
$ws.Range("Q9").Value = 205.5355265724078
$ws.Range("R9").Value = 1849.81973915167
$ws.Range("S9").Value = 0.2494498943116692
$ws.Range("T9").Value = 0.2494498943116693
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 64.69033266666666
$ws.Range("H10").Value = 194.070998
$ws.Range("I10").Value = 0.4264136233238083
$ws.Range("J10").Value = 0.4264136233238083
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.384523666666666
$ws.Range("N10").Value = 4.153570999999999
$ws.Range("O10").Value = 0.2549207111021641
$ws.Range("P10").Value = 0.2549207111021641
$ws.Range("Q10").Value = 89.56529658153976
$ws.Range("R10").Value = 806.0876692338579
$ws.Range("S10").Value = 0.1087016640813556
$ws.Range("T10").Value = 0.1087016640813556